# Apply the edits described by the commit "agregue los archivos de pruebas"
# to the "Pruebas" worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pruebas")
$ws.Activate()

# --- Header cells: expand short labels into full text ---
$ws.Range("A1").Value = "Proyecto BlueWeb"
$ws.Range("A2").Value = "Pantalla Ciudad"

# --- Re-capitalize a handful of existing "resultado"/"funcionalidad" cells ---
$ws.Range("F4").Value = "Grafica con datos mostrados con exito."
$ws.Range("F5").Value = "Elementos eliminados con exito "
$ws.Range("E6").Value = "Agrega registos a la tabla al llenar los campos solicitados"
$ws.Range("F6").Value = "Elementos agregados con éxito"
$ws.Range("E7").Value = "Modifica los campos de la tabla"
$ws.Range("F7").Value = "Elementos modificados exitosamente"
$ws.Range("B8").Value = "Editar solo un dato de la fila"
$ws.Range("E8").Value = "Edita un campo y manda ese valor a la funcion de editar, ignorando los otros"
$ws.Range("F8").Value = "Edita un solo campo dejando los otros iguales"
$ws.Range("E9").Value = "Edita 2 campos del registro dejando el tercero igual"
$ws.Range("F9").Value = "Edita solo los 2 campos exitosamente"

# --- Fill in the previously-empty row 10 with a new test case (#7) ---
# Copy the date format from D9 first so D10 picks up the same number format
# (m/d/yyyy) style instead of Excel inventing a brand-new style entry.
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Presionar boton enter"
$ws.Range("C10").Value = "al presionar el boton enter al agregar se enviar el formulario correctamente"
$ws.Range("D10").Value = 44431
$ws.Range("E10").Value = "Cuando se abre el modal para eliminar, y se llenan lo campos, al presionar enter se activa el boton de agregar"
$ws.Range("F10").Value = "El boton de agregar, agrega correctamente al presionar enter"

$ws.Rows.Item(10).RowHeight = 43.5

# --- Move the active selection to match the post-edit cursor position ---
$ws.Range("F9").Select()
